$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 5.675800000000002
$ws.Range("C3").Value = -12.20399999999999
$ws.Range("E3").Value = 16.42790000000001
$ws.Range("E6").Value = 16.40330000000001
$ws.Range("D8").Value = -9.048099999999998
$ws.Range("D11").Value = -7.5583
$ws.Range("A12").Value = -21.57320000000001
$ws.Range("B14").Value = 6.029899999999997
$ws.Range("D14").Value = -7.590700000000003
$ws.Range("D15").Value = -8.316299999999996
$ws.Range("B26").Value = 4.001000000000001
$ws.Range("E27").Value = 16.6332
$ws.Range("C30").Value = -13.6955
$ws.Range("B31").Value = 4.570900000000002
$ws.Range("A32").Value = -21.25170000000001
$ws.Range("E33").Value = 17.16170000000002
$ws.Range("B35").Value = 9.3452
$ws.Range("A36").Value = -19.425
$ws.Range("D36").Value = -7.892900000000003
$ws.Range("B37").Value = 9.011700000000001
$ws.Range("A38").Value = -19.2708
$ws.Range("E39").Value = 15.7743
$ws.Range("C44").Value = -13.55499999999999
$ws.Range("B45").Value = 5.447299999999999
$ws.Range("A46").Value = -21.8042
$ws.Range("E47").Value = 16.6244
$ws.Range("A54").Value = -22.2843
$ws.Range("E54").Value = 16.91679999999999
$ws.Range("A55").Value = -21.93489999999998
$ws.Range("E56").Value = 16.16690000000001
$ws.Range("B57").Value = 4.828299999999997
$ws.Range("C58").Value = -12.9883
$ws.Range("E58").Value = 16.6306
$ws.Range("D64").Value = -7.550699999999996
$ws.Range("E66").Value = 17.13650000000001
$ws.Range("A67").Value = -21.41709999999998
$ws.Range("A69").Value = -21.56429999999998
$ws.Range("A72").Value = -21.91310000000001
$ws.Range("E72").Value = 16.6883
$ws.Range("E82").Value = 16.59840000000001
$ws.Range("E83").Value = 16.87939999999999
$ws.Range("C84").Value = -13.87409999999999
$ws.Range("C89").Value = -11.0496
$ws.Range("D89").Value = -5.830300000000002
$ws.Range("A91").Value = -21.5674
$ws.Range("C91").Value = -11.0867
$ws.Range("C92").Value = -11.3661
$ws.Range("A99").Value = -20.17529999999998
$ws.Range("B100").Value = 4.815299999999997
$ws.Range("B102").Value = 7.9859
$ws.Range("C102").Value = -13.2627
